# Risk assessment spreadsheet updated as per meeting 8.1.19
#
# Changes:
#  1. "Over Scoping" row (row 20): Likelihood (D) and Overall Risk Level (F)
#     downgraded from High to Medium.
#  2. "Code Understanding" row (row 21): Likelihood (D) upgraded from Low to
#     Medium.
#  3. New risk added in row 26: "Tutorial (intuitiveness)" with
#     Likelihood=Medium, Severity=High, Overall Risk Level=High.
#  4. Selection left on D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: "Over Scoping" -> Likelihood & Overall Risk Level: High -> Medium
$ws.Range("D17").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("D20").Value = "Medium"
$ws.Range("F20").Value = "Medium"

# --- Row 21: "Code Understanding" -> Likelihood: Low -> Medium
$ws.Range("D17").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D21").Value = "Medium"

# --- Row 26: new risk "Tutorial (intuitiveness)"
$ws.Range("A26").Value = "Tutorial (intuitiveness)"

$ws.Range("D17").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = "Medium"

$ws.Range("E15").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = "High"

$ws.Range("E15").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Font.Name = "Calibri"
$ws.Range("F26").Value = "High"

# --- Leave the cursor where the author left it
$ws.Range("D9").Select()
